$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the player it describes
$ws.Name = "Prasidh Krishna"

# Insert a new column before column A for "matchNo" -> shifts existing
# teamName..result columns from A:L to B:M
$ws.Columns.Item(1).Insert()

# Insert a new row before row 2 for the "15th" match entry -> shifts the
# existing "18th" (Rajasthan Royals) row from row 2 down to row 3
$ws.Rows.Item(2).Insert()

# Force text formatting on the numeric-looking stat columns so values such
# as "0" / "1" / "0.00" are kept as text instead of being re-interpreted
# as numbers (matching the source data, which stores everything as text).
$textCells = @("E2", "F2", "G2", "H2", "I2", "E3", "F3", "G3", "H3")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Header row
$ws.Range("A1").Value2 = "matchNo"
$ws.Range("B1").Value2 = "teamName"
$ws.Range("C1").Value2 = "batterName"
$ws.Range("D1").Value2 = "states"
$ws.Range("E1").Value2 = "runs"
$ws.Range("F1").Value2 = "balls"
$ws.Range("G1").Value2 = "fours"
$ws.Range("H1").Value2 = "sixes"
$ws.Range("I1").Value2 = "sr"
$ws.Range("J1").Value2 = "opponentTeamName"
$ws.Range("K1").Value2 = "venue"
$ws.Range("L1").Value2 = "date"
$ws.Range("M1").Value2 = "result"

# Row 2 - new "15th" match entry (vs Chennai Super Kings)
$ws.Range("A2").Value2 = "15th"
$ws.Range("B2").Value2 = "Kolkata Knight Riders"
$ws.Range("C2").Value2 = "Prasidh Krishna"
$ws.Range("D2").Value2 = "run out (Chahar/Thakur)"
$ws.Range("E2").Value2 = "0"
$ws.Range("F2").Value2 = "1"
$ws.Range("G2").Value2 = "0"
$ws.Range("H2").Value2 = "0"
$ws.Range("I2").Value2 = "0.00"
$ws.Range("J2").Value2 = "Chennai Super Kings"
$ws.Range("K2").Value2 = "Wankhede"
$ws.Range("L2").Value2 = "April 21"
$ws.Range("M2").Value2 = "Super Kings won by 18 runs"

# Row 3 - existing "18th" match entry, now shifted down from row 2 and
# right one column (A:L -> B:M); only the new "matchNo" cell (A3) needs to
# be written, the rest of the row already holds the correct shifted data.
$ws.Range("A3").Value2 = "18th"

Write-Host "Edit complete"
